$d = $word.ActiveDocument

# "Chỉnh lại mẫu 26": remove the placeholder merge-field run
# "vnpt.SiteAddress" that follows "Địa chỉ: " in the "Bên A" address line,
# leaving just the label run "Địa chỉ: " in place.
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("vnpt.SiteAddress", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

if (-not $found) {
    throw "Could not find 'vnpt.SiteAddress' placeholder text to remove."
}
